$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Insert a new row above the old row 13 (LATENCY_WINDOW_SAMPLES), pushing the
# remaining parameter rows down by one.
$ws.Rows.Item(13).Insert()

# Populate the new row with the "FB_JERK_NEG_THD" parameter.
# Write the Description column (E) before the Parameter column (A) so the
# new shared-string entries land in the same order as the authoritative
# edit (description string first, then the parameter-name string).
$ws.Cells.Item(13, 5).Value2 = "Negative jerk threshold to detect the start of full braking"
$ws.Cells.Item(13, 1).Value2 = "FB_JERK_NEG_THD"
$ws.Cells.Item(13, 2).Value2 = -20
$ws.Cells.Item(13, 3).Value2 = "float"
$ws.Cells.Item(13, 4).Value2 = "m/s³"
$ws.Cells.Item(13, 6).Value2 = "AebKpiExtractor"

# Switch the active sheet/selection to the params sheet at the new row.
$ws.Activate()
[void]$ws.Range("A13").Select()
